# Complete Missing Companies Sorter
# Update the "Updated" date column (B2:B61) from 2023-07-14 (45121) to 2023-07-18 (45125)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldDate = 45121
$newDate = 45125

$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 61 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    if ($cell.Value2 -eq $oldDate) {
        $cell.Value = $newDate
    }
}
